# Update column F ("dSF") values for specific rows to reflect the
# repulled/recalculated data, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 0
    9  = 0
    19 = 0
    20 = -3
    21 = 4
    22 = -3
    26 = -1
    27 = -2
    30 = -2
    31 = -5
    35 = -1
    36 = -2
    39 = 0
    47 = -4
    48 = -1
    51 = -1
    53 = -3
    56 = 2
    60 = -1
    62 = -2
    63 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
